$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6900.75
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 8867.666999999999
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 8867.666999999999
$ws.Range("M9").Value = -831
$ws.Range("N9").Value = -9205.666999999999
$ws.Range("H19").Value = 1756.9048
$ws.Range("I19").Value = 2797
$ws.Range("K19").Value = 2797
$ws.Range("M19").Value = -2622
$ws.Range("H116").Value = 4294.15
$ws.Range("I116").Value = 3724.5833
$ws.Range("K116").Value = 3724.5833
$ws.Range("M116").Value = -282.5832999999998
$ws.Range("H132").Value = 1421.4814
$ws.Range("I132").Value = 1103.125
$ws.Range("K132").Value = 3309.375
$ws.Range("M132").Value = -779.375
$ws.Range("H137").Value = 2779.9546
$ws.Range("I137").Value = 1808.5238
$ws.Range("J137").Value = 3233.2888
$ws.Range("K137").Value = 5425.5714
$ws.Range("L137").Value = 9699.866399999999
$ws.Range("M137").Value = -2875.5714
$ws.Range("N137").Value = -14799.8664
$ws.Range("H138").Value = 2462.3872
$ws.Range("I138").Value = 1022.7778
$ws.Range("J138").Value = 3371.614
$ws.Range("K138").Value = 3068.3334
$ws.Range("L138").Value = 10114.842
$ws.Range("M138").Value = 2071.6666
$ws.Range("N138").Value = -20394.842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9742.546
$ws.Range("I2").Value = 458.5
$ws.Range("K2").Value = 458.5
$ws.Range("M2").Value = -345.5
$ws.Range("H3").Value = 15851.25
$ws.Range("I3").Value = 4702
$ws.Range("J3").Value = 27000.5
$ws.Range("K3").Value = 4702
$ws.Range("L3").Value = 27000.5
$ws.Range("M3").Value = -4587
$ws.Range("N3").Value = -27230.5
$ws.Range("H22").Value = 3016.2
$ws.Range("I22").Value = 3520.25
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3520.25
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -3221.25
$ws.Range("N22").Value = -1598
$ws.Range("H32").Value = 4171.1333
$ws.Range("I32").Value = 3183.6191
$ws.Range("K32").Value = 3183.6191
$ws.Range("M32").Value = -2896.6191
$ws.Range("H45").Value = 3268.9
$ws.Range("I45").Value = 2961.125
$ws.Range("J45").Value = 4500
$ws.Range("K45").Value = 2961.125
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = -2584.125
$ws.Range("N45").Value = -5254
$ws.Range("H61").Value = 4292.5293
$ws.Range("I61").Value = 3543.394
$ws.Range("K61").Value = 3543.394
$ws.Range("M61").Value = -3331.394
$ws.Range("H116").Value = 9742.546
$ws.Range("I116").Value = 458.5
$ws.Range("K116").Value = 458.5
$ws.Range("M116").Value = 1835.5
$ws.Range("H132").Value = 2445.353
$ws.Range("I132").Value = 1810.742
$ws.Range("K132").Value = 5432.226
$ws.Range("M132").Value = -2902.226
$ws.Range("H136").Value = 4292.5293
$ws.Range("I136").Value = 3543.394
$ws.Range("K136").Value = 10630.182
$ws.Range("M136").Value = -8080.181999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9742.546
$ws.Range("I3").Value = 458.5
$ws.Range("K3").Value = 458.5
$ws.Range("M3").Value = -344.5
$ws.Range("H8").Value = 25001
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("H20").Value = 4606.524
$ws.Range("I20").Value = 4141.5625
$ws.Range("K20").Value = 4141.5625
$ws.Range("M20").Value = -3894.5625
$ws.Range("H134").Value = 3389
$ws.Range("I134").Value = 2902.2856
$ws.Range("K134").Value = 8706.856800000001
$ws.Range("M134").Value = -6171.856800000001
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7000.5713
$ws.Range("I58").Value = 2111
$ws.Range("J58").Value = 15801.8
$ws.Range("K58").Value = 2111
$ws.Range("L58").Value = 15801.8
$ws.Range("M58").Value = -1908
$ws.Range("N58").Value = -16207.8
$ws.Range("H99").Value = 3956.5715
$ws.Range("I99").Value = 3519.4
$ws.Range("J99").Value = 5049.5
$ws.Range("K99").Value = 3519.4
$ws.Range("L99").Value = 5049.5
$ws.Range("M99").Value = -2021.4
$ws.Range("N99").Value = -8045.5
$ws.Range("H126").Value = 3956.5715
$ws.Range("I126").Value = 3519.4
$ws.Range("J126").Value = 5049.5
$ws.Range("K126").Value = 10558.2
$ws.Range("L126").Value = 15148.5
$ws.Range("M126").Value = -8088.200000000001
$ws.Range("N126").Value = -20088.5
$ws.Range("H132").Value = 2272.5112
$ws.Range("I132").Value = 1521.2059
$ws.Range("J132").Value = 4594.727
$ws.Range("K132").Value = 4563.6177
$ws.Range("L132").Value = 13784.181
$ws.Range("M132").Value = -2033.6177
$ws.Range("N132").Value = -18844.181
$ws.Range("H134").Value = 2151.926
$ws.Range("I134").Value = 1647.826
$ws.Range("K134").Value = 4943.478
$ws.Range("M134").Value = -2408.478
$ws.Range("H136").Value = 7000.5713
$ws.Range("I136").Value = 2111
$ws.Range("J136").Value = 15801.8
$ws.Range("K136").Value = 6333
$ws.Range("L136").Value = 47405.39999999999
$ws.Range("M136").Value = -3783
$ws.Range("N136").Value = -52505.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 260889.67
$ws.Range("J37").Value = 260889.67
$ws.Range("L37").Value = 782669.01
$ws.Range("N37").Value = -782893.01
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H63").Value = 16890.334
$ws.Range("I63").Value = 3999
$ws.Range("J63").Value = 18501.75
$ws.Range("K63").Value = 11997
$ws.Range("L63").Value = 55505.25
$ws.Range("M63").Value = -11248
$ws.Range("N63").Value = -57003.25
$ws.Range("H64").Value = 58828556
$ws.Range("I64").Value = 250001070
$ws.Range("J64").Value = 6241.6924
$ws.Range("K64").Value = 750003210
$ws.Range("L64").Value = 18725.0772
$ws.Range("M64").Value = -750002940
$ws.Range("N64").Value = -19265.0772
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H66").Value = 16890.334
$ws.Range("I66").Value = 3999
$ws.Range("J66").Value = 18501.75
$ws.Range("K66").Value = 35991
$ws.Range("L66").Value = 166515.75
$ws.Range("M66").Value = -32247
$ws.Range("N66").Value = -174003.75
$ws.Range("H67").Value = 58828556
$ws.Range("I67").Value = 250001070
$ws.Range("J67").Value = 6241.6924
$ws.Range("K67").Value = 750003210
$ws.Range("L67").Value = 18725.0772
$ws.Range("M67").Value = -750002274
$ws.Range("N67").Value = -20597.0772
$ws.Range("H69").Value = 9402.200000000001
$ws.Range("J69").Value = 9402.200000000001
$ws.Range("L69").Value = 28206.6
$ws.Range("N69").Value = -29828.6
$ws.Range("H70").Value = 8666.333000000001
$ws.Range("I70").Value = 5499.5
$ws.Range("K70").Value = 16498.5
$ws.Range("M70").Value = -16183.5
$ws.Range("H72").Value = 9402.200000000001
$ws.Range("J72").Value = 9402.200000000001
$ws.Range("L72").Value = 84619.8
$ws.Range("N72").Value = -92731.8
$ws.Range("H73").Value = 8666.333000000001
$ws.Range("I73").Value = 5499.5
$ws.Range("K73").Value = 16498.5
$ws.Range("M73").Value = -15406.5
$ws.Range("H137").Value = 77558.30499999999
$ws.Range("I137").Value = 688.0833
$ws.Range("K137").Value = 2064.2499
$ws.Range("M137").Value = 3035.7501
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2730.0322
$ws.Range("I102").Value = 2173.238
$ws.Range("K102").Value = 2173.238
$ws.Range("M102").Value = -551.2379999999998
$ws.Range("H132").Value = 2273.5
$ws.Range("I132").Value = 1254.6111
$ws.Range("J132").Value = 4566
$ws.Range("K132").Value = 3763.8333
$ws.Range("L132").Value = 13698
$ws.Range("M132").Value = -1233.8333
$ws.Range("N132").Value = -18758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9603.440000000001
$ws.Range("I136").Value = 7353.8
$ws.Range("K136").Value = 22061.4
$ws.Range("M136").Value = -19511.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2716.1667
$ws.Range("I122").Value = 1719.7273
$ws.Range("K122").Value = 5159.1819
$ws.Range("M122").Value = -2709.1819
